# This script reproduces the "chore: update Sheets via scheduled runner" commit.
# It rewrites a set of market-price-derived cells (currentAveragePrice* /
# LevePrice* / LeveProfit* columns H,I,J,K,L,M,N) on several sheets of the
# Jenova_Profits workbook to their refreshed values, leaving every other
# cell (labels, Leve metadata, formulas, etc.) untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1069.2
$ws.Range("I15").Value = 1069.2
$ws.Range("K15").Value = 3207.6
$ws.Range("M15").Value = -3038.6

# Row 28
$ws.Range("H28").Value = 102338.8
$ws.Range("I28").Value = 102338.8
$ws.Range("K28").Value = 102338.8
$ws.Range("M28").Value = -101853.8

# Row 40
$ws.Range("H40").Value = 6920.1333
$ws.Range("I40").Value = 4978.1113
$ws.Range("J40").Value = 9833.166999999999
$ws.Range("K40").Value = 4978.1113
$ws.Range("L40").Value = 9833.166999999999
$ws.Range("M40").Value = -4803.1113
$ws.Range("N40").Value = -10183.167

# Row 43
$ws.Range("H43").Value = 1846.6
$ws.Range("I43").Value = 1866.6666
$ws.Range("J43").Value = 1816.5
$ws.Range("K43").Value = 1866.6666
$ws.Range("L43").Value = 1816.5
$ws.Range("M43").Value = -1797.6666
$ws.Range("N43").Value = -1954.5

# Row 64
$ws.Range("H64").Value = 6625
$ws.Range("J64").Value = 6625
$ws.Range("L64").Value = 6625
$ws.Range("N64").Value = -7121

# Row 67
$ws.Range("H67").Value = 6625
$ws.Range("J67").Value = 6625
$ws.Range("L67").Value = 6625
$ws.Range("N67").Value = -8341

# Row 76
$ws.Range("H76").Value = 111117330
$ws.Range("I76").Value = 5319.3335
$ws.Range("J76").Value = 166673330
$ws.Range("K76").Value = 5319.3335
$ws.Range("L76").Value = 166673330
$ws.Range("M76").Value = -5004.3335
$ws.Range("N76").Value = -166673960

# Row 79
$ws.Range("H79").Value = 111117330
$ws.Range("I79").Value = 5319.3335
$ws.Range("J79").Value = 166673330
$ws.Range("K79").Value = 5319.3335
$ws.Range("L79").Value = 166673330
$ws.Range("M79").Value = -4227.3335
$ws.Range("N79").Value = -166675514

# Row 86
$ws.Range("H86").Value = 2291473.5
$ws.Range("I86").Value = 2768.2222
$ws.Range("K86").Value = 2768.2222
$ws.Range("M86").Value = -1645.2222

# Row 89
$ws.Range("H89").Value = 2291473.5
$ws.Range("I89").Value = 2768.2222
$ws.Range("K89").Value = 13841.111
$ws.Range("M89").Value = -8225.111000000001

# Row 106
$ws.Range("H106").Value = 2567.2222
$ws.Range("I106").Value = 2424.4119
$ws.Range("K106").Value = 2424.4119
$ws.Range("M106").Value = -1793.4119

# Row 138
$ws.Range("H138").Value = 5150.4346
$ws.Range("I138").Value = 2443.3333
$ws.Range("J138").Value = 5808.919
$ws.Range("K138").Value = 7329.999899999999
$ws.Range("L138").Value = 17426.757
$ws.Range("M138").Value = -2189.999899999999
$ws.Range("N138").Value = -27706.757

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 642.2381
$ws.Range("I2").Value = 640.7368
$ws.Range("K2").Value = 640.7368
$ws.Range("M2").Value = -527.7368

# Row 32
$ws.Range("H32").Value = 1572.21
$ws.Range("I32").Value = 1542.6364
$ws.Range("K32").Value = 1542.6364
$ws.Range("M32").Value = -1255.6364

# Row 45
$ws.Range("H45").Value = 2366.3076
$ws.Range("I45").Value = 1695.7778
$ws.Range("J45").Value = 3875
$ws.Range("K45").Value = 1695.7778
$ws.Range("L45").Value = 3875
$ws.Range("M45").Value = -1318.7778
$ws.Range("N45").Value = -4629

# Row 102
$ws.Range("H102").Value = 1902
$ws.Range("I102").Value = 1442.2
$ws.Range("K102").Value = 1442.2
$ws.Range("M102").Value = 179.8

# Row 116
$ws.Range("H116").Value = 642.2381
$ws.Range("I116").Value = 640.7368
$ws.Range("K116").Value = 640.7368
$ws.Range("M116").Value = 1653.2632

# Row 132
$ws.Range("H132").Value = 5523.5713
$ws.Range("I132").Value = 5829.4375
$ws.Range("J132").Value = 4544.8
$ws.Range("K132").Value = 17488.3125
$ws.Range("L132").Value = 13634.4
$ws.Range("M132").Value = -14958.3125
$ws.Range("N132").Value = -18694.4

# Row 134
$ws.Range("H134").Value = 50000
$ws.Range("I134").Value = 50000
$ws.Range("K134").Value = 50000
$ws.Range("M134").Value = -44930

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 642.2381
$ws.Range("I3").Value = 640.7368
$ws.Range("K3").Value = 640.7368
$ws.Range("M3").Value = -526.7368

# Row 134
$ws.Range("H134").Value = 28250.047
$ws.Range("I134").Value = 4040.6487
$ws.Range("K134").Value = 12121.9461
$ws.Range("M134").Value = -9586.946100000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 45737.72
$ws.Range("J31").Value = 65528.35
$ws.Range("L31").Value = 65528.35
$ws.Range("N31").Value = -66118.35000000001

# Row 34
$ws.Range("H34").Value = 45737.72
$ws.Range("J34").Value = 65528.35
$ws.Range("L34").Value = 65528.35
$ws.Range("N34").Value = -65932.35000000001

# Row 94
$ws.Range("H94").Value = 977.7692
$ws.Range("I94").Value = 604
$ws.Range("J94").Value = 1089.9
$ws.Range("K94").Value = 604
$ws.Range("L94").Value = 1089.9
$ws.Range("M94").Value = -153
$ws.Range("N94").Value = -1991.9

# Row 99
$ws.Range("H99").Value = 4047.3
$ws.Range("I99").Value = 3745.8333
$ws.Range("K99").Value = 3745.8333
$ws.Range("M99").Value = -2247.8333

# Row 105
$ws.Range("H105").Value = 742.3333
$ws.Range("I105").Value = 662.53845
$ws.Range("K105").Value = 662.53845
$ws.Range("M105").Value = 1084.46155

# Row 126
$ws.Range("H126").Value = 4047.3
$ws.Range("I126").Value = 3745.8333
$ws.Range("K126").Value = 11237.4999
$ws.Range("M126").Value = -8767.499899999999

# Row 132
$ws.Range("H132").Value = 2377.1428
$ws.Range("J132").Value = 3145.1428
$ws.Range("L132").Value = 9435.428400000001
$ws.Range("N132").Value = -14495.4284

$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 121973.75
$ws.Range("J37").Value = 121973.75
$ws.Range("L37").Value = 365921.25
$ws.Range("N37").Value = -366145.25

# Row 117
$ws.Range("H117").Value = 872
$ws.Range("I117").Value = 674.75
$ws.Range("J117").Value = 950.9
$ws.Range("K117").Value = 2024.25
$ws.Range("L117").Value = 2852.7
$ws.Range("M117").Value = 1417.75
$ws.Range("N117").Value = -9736.700000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2920.879
$ws.Range("I40").Value = 2960.0312
$ws.Range("K40").Value = 2960.0312
$ws.Range("M40").Value = -2824.0312

# Row 122
$ws.Range("H122").Value = 4665.1665
$ws.Range("I122").Value = 4398.4
$ws.Range("K122").Value = 13195.2
$ws.Range("M122").Value = -10745.2

# Row 132
$ws.Range("H132").Value = 2951.6667
$ws.Range("I132").Value = 1142.1
$ws.Range("K132").Value = 3426.3
$ws.Range("M132").Value = -896.2999999999997

# Row 136
$ws.Range("H136").Value = 364090.22
$ws.Range("I136").Value = 719143.6
$ws.Range("K136").Value = 2157430.8
$ws.Range("M136").Value = -2154880.8

$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 60000
$ws.Range("I51").Value = 60000
$ws.Range("K51").Value = 60000
$ws.Range("M51").Value = -59490

# Row 52
$ws.Range("H52").Value = 60000
$ws.Range("I52").Value = 60000
$ws.Range("K52").Value = 60000
$ws.Range("M52").Value = -59774
